$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update 想去人数 (want-to-go count) values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 1050
$wsExpo.Range("F3").Value = 312
$wsExpo.Range("F4").Value = 2833

# Sheet "全部类型" (all types) - same rows duplicated, update the same values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1050
$wsAll.Range("F5").Value = 312
$wsAll.Range("F6").Value = 2833
